$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-02"

# Update the header label cell (I1) with the new "through" date text
$ws.Range("I1").Value = "2022 (through 03-02)"

# Update the March (row 4) total for the latest column (I) from 8 to 9
$ws.Range("I4").Value = 9

# Update the yearly Total row (row 14) for the latest column (I) from 308 to 309
$ws.Range("I14").Value = 309
